$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last-modified timestamp string (P1) in row 1
$ws.Range("P1").Value = "2018-09-03 22:04:13"

# Add two new pairs of columns (Q1:T1) mirroring the existing "campo"/0 pattern
$ws.Range("Q1").Value = "campo"
$ws.Range("R1").Value = 0
$ws.Range("S1").Value = "campo"
$ws.Range("T1").Value = 0
